# The workbook originally ships with a single sheet ("sheet0") holding a
# handful of merge-field placeholders. This change makes every sheet in
# the workbook render with that same data (commit: "render all
# excel-sheets with thee same data") by duplicating "sheet0" into a new
# "sheet1", and leaves the duplicate as the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing sheet right after itself - this brings along the
# cell values (and their shared-string usage) as well as the existing
# per-cell styles in one step.
$ws.Copy($null, $ws)

$ws2 = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -ne $ws.Name) {
        $ws2 = $sheet
    }
}
$ws2.Name = "sheet1"

# The freshly duplicated sheet becomes the active tab/sheet.
$ws2.Select()

# Refresh the on-sheet selection from the last data row to the second row
# on both sheets.
$ws.Range("A2").Select()
$ws2.Range("A2").Select()

# Give the template/body rows (every row below the row-1 header) their own
# distinct cell style on both sheets, rather than sharing the sheet's
# default style.
$bodyRows = @(2, 4, 5)
foreach ($r in $bodyRows) {
    $ws.Cells.Item($r, 1).HorizontalAlignment = -4131
    $ws2.Cells.Item($r, 1).HorizontalAlignment = -4131
}

$ws2.Select()
